# Update "想去人数" (column F) figures and a couple of "最低票价" (column G)
# sold-out markers to match the refreshed data snapshot.
#
# Column letters used below: F = 6, G = 7

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览 (Exhibition)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(2, 6).Value2  = 12611
$ws1.Cells.Item(3, 6).Value2  = 7045
$ws1.Cells.Item(10, 6).Value2 = 989
$ws1.Cells.Item(10, 7).Value2 = "已售罄"
$ws1.Cells.Item(12, 6).Value2 = 343
$ws1.Cells.Item(13, 6).Value2 = 991
$ws1.Cells.Item(15, 6).Value2 = 1011
$ws1.Cells.Item(16, 6).Value2 = 510
$ws1.Cells.Item(20, 6).Value2 = 269
$ws1.Cells.Item(22, 6).Value2 = 43
$ws1.Cells.Item(23, 6).Value2 = 129
$ws1.Cells.Item(25, 6).Value2 = 5193
$ws1.Cells.Item(27, 6).Value2 = 1397
$ws1.Cells.Item(28, 6).Value2 = 299
$ws1.Cells.Item(29, 6).Value2 = 1248
$ws1.Cells.Item(30, 6).Value2 = 1248
$ws1.Cells.Item(31, 6).Value2 = 1323
$ws1.Cells.Item(32, 6).Value2 = 1
$ws1.Cells.Item(35, 6).Value2 = 3719

# ---------------------------------------------------------------------------
# Sheet 2: 演出 (Performance)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(8, 6).Value2 = 37

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活 (Local life)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Cells.Item(2, 6).Value2 = 9238
$ws3.Cells.Item(3, 6).Value2 = 551
$ws3.Cells.Item(4, 6).Value2 = 1968

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 (All types)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Cells.Item(2, 6).Value2  = 9238
$ws4.Cells.Item(3, 6).Value2  = 551
$ws4.Cells.Item(4, 6).Value2  = 1968
$ws4.Cells.Item(6, 6).Value2  = 12611
$ws4.Cells.Item(7, 6).Value2  = 7045
$ws4.Cells.Item(13, 6).Value2 = 989
$ws4.Cells.Item(13, 7).Value2 = "已售罄"
$ws4.Cells.Item(15, 6).Value2 = 343
$ws4.Cells.Item(16, 6).Value2 = 991
$ws4.Cells.Item(18, 6).Value2 = 1011
$ws4.Cells.Item(19, 6).Value2 = 510
$ws4.Cells.Item(23, 6).Value2 = 269
$ws4.Cells.Item(25, 6).Value2 = 43
$ws4.Cells.Item(28, 6).Value2 = 37
$ws4.Cells.Item(32, 6).Value2 = 5193
$ws4.Cells.Item(34, 6).Value2 = 1397
$ws4.Cells.Item(37, 6).Value2 = 299
$ws4.Cells.Item(39, 6).Value2 = 1248
$ws4.Cells.Item(40, 6).Value2 = 1248
$ws4.Cells.Item(41, 6).Value2 = 1323
$ws4.Cells.Item(42, 6).Value2 = 1
$ws4.Cells.Item(48, 6).Value2 = 3719
